# Update countries & provincias Spain
# Applies the COVID-19 data refresh (23 May 2020, 14:05 -> 14:35) to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 14:35"

# --- Row 6: Brasil (unchanged name, updated figures) -------------------------
$ws.Cells.Item(6, 2).Value = 333937   # Casos totales
$ws.Cells.Item(6, 3).Value = 3047     # Nuevos casos
$ws.Cells.Item(6, 5).Value = 177362   # Recuperados
$ws.Cells.Item(6, 7).Value = 97       # Muertes hoy
$ws.Cells.Item(6, 8).Value = 21145    # Muertes

# --- Row 14: India -------------------------------------------------------
$ws.Cells.Item(14, 2).Value = 126405
$ws.Cells.Item(14, 3).Value = 1611
$ws.Cells.Item(14, 4).Value = 52272
$ws.Cells.Item(14, 5).Value = 70379

# --- Row 23: Paises Bajos -------------------------------------------------
$ws.Cells.Item(23, 2).Value = 45064
$ws.Cells.Item(23, 3).Value = 176
$ws.Cells.Item(23, 7).Value = 23
$ws.Cells.Item(23, 8).Value = 5811

# --- Row 31: Portugal ------------------------------------------------------
$ws.Cells.Item(31, 2).Value = 30471
$ws.Cells.Item(31, 3).Value = 271
$ws.Cells.Item(31, 4).Value = 7705
$ws.Cells.Item(31, 5).Value = 21464
$ws.Cells.Item(31, 7).Value = 13
$ws.Cells.Item(31, 8).Value = 1302

# --- Row 63: Moldavia -------------------------------------------------------
$ws.Cells.Item(63, 4).Value = 3452
$ws.Cells.Item(63, 5).Value = 3157
$ws.Cells.Item(63, 7).Value = 1
$ws.Cells.Item(63, 8).Value = 238

# --- Rows 75-77: Uzbekistan overtakes Guinea & Tailandia in the ranking ------
# New order: Uzbekistan (updated figures), Guinea (old row75 figures),
# Tailandia (old row76 figures). Senegal (row78) is unaffected.
$ws.Cells.Item(75, 1).Value = "Uzbekistan"
$ws.Cells.Item(75, 2).Value = 3089
$ws.Cells.Item(75, 3).Value = 61
$ws.Cells.Item(75, 4).Value = 2525
$ws.Cells.Item(75, 5).Value = 551
$ws.Cells.Item(75, 8).Value = 13

$ws.Cells.Item(76, 1).Value = "Guinea"
$ws.Cells.Item(76, 2).Value = 3067
$ws.Cells.Item(76, 3).Value = 0
$ws.Cells.Item(76, 4).Value = 1575
$ws.Cells.Item(76, 5).Value = 1473
$ws.Cells.Item(76, 8).Value = 19

$ws.Cells.Item(77, 1).Value = "Tailandia"
$ws.Cells.Item(77, 2).Value = 3040
$ws.Cells.Item(77, 3).Value = 3
$ws.Cells.Item(77, 4).Value = 2916
$ws.Cells.Item(77, 5).Value = 68
$ws.Cells.Item(77, 8).Value = 56

# --- Row 86: Croacia ---------------------------------------------------------
$ws.Cells.Item(86, 4).Value = 2023
$ws.Cells.Item(86, 5).Value = 121

# --- Row 105: Sri Lanka -------------------------------------------------------
$ws.Cells.Item(105, 2).Value = 1085
$ws.Cells.Item(105, 3).Value = 17
$ws.Cells.Item(105, 5).Value = 416

# --- Rows 198-200: Santa Lucia / Belice / Nueva Caledonia rotate -------------
# Old: 198=Santa Lucia, 199=Belice, 200=Nueva Caledonia
# New: 198=Belice, 199=Nueva Caledonia, 200=Santa Lucia
$ws.Cells.Item(198, 1).Value = "Belice"
$ws.Cells.Item(198, 4).Value = 16
$ws.Cells.Item(198, 8).Value = 2

$ws.Cells.Item(199, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(199, 4).Value = 18
$ws.Cells.Item(199, 8).Value = 0

$ws.Cells.Item(200, 1).Value = "Santa Lucia"
$ws.Cells.Item(200, 4).Value = 18
$ws.Cells.Item(200, 8).Value = 0

# --- Rows 210-211: Montserrat / Groenlandia swap ------------------------------
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = "Groenlandia"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# --- Rows 215-216: San Bartolome / Bonaire, San Eustaquio y Saba swap ---------
# (All values tied at 6,0,6,0,0,0,0 for both countries, so only the names move.)
# Old: 214=Sahara Occidental (unaffected), 215=San Bartolome, 216=Bonaire...
# New: 214=Sahara Occidental (unaffected), 215=Bonaire..., 216=San Bartolome
$ws.Cells.Item(215, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(216, 1).Value = "San Bartolome"
